$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: LeftImage -> leftImgPath, RightImage -> rightImgPath
$ws.Range("B1").Value = "leftImgPath"
$ws.Range("C1").Value = "rightImgPath"

# Reflect final selected cell as saved in the workbook
$ws.Range("H12").Select()
